# Adjusted to Naming Convention
# Renames the Pset fuel-type codes (EPPxxx -> P*xxx*), moves the label
# column from D to F (new style), and adds three new fuel rows
# (Peat, Distillate, HFO) while re-ordering the renewable rows so that
# Geo/Ocean swap places - matching the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Column widths / view tweaks
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.6328125
$ws.Columns.Item(4).ColumnWidth = 13.453125
$ws.Columns.Item(7).ColumnWidth = 10.7265625
$ws.Columns.Item(8).ColumnWidth = 12.453125
$ws.Columns.Item(10).ColumnWidth = 11.08984375

# ---------------------------------------------------------------
# 2. Clear the old column-D labels for the fuel rows (9-17); the
#    label now lives in column F.
# ---------------------------------------------------------------
$ws.Range("D9:D17").ClearContents()

# ---------------------------------------------------------------
# 3. Write the new fuel-row table (rows 9-20).
#    Column F = fuel code label (styled with the CJK font),
#    Column A = RNW flag (1 for renewables, blank for fossil fuels),
#    Columns J-M = fixed lookups, N-T = the IF() formulas.
# ---------------------------------------------------------------
$fuelRows = @(
    @{ Row = 9;  Label = "P*GAS*"; Flag = $false },
    @{ Row = 10; Label = "P*COA*"; Flag = $false },
    @{ Row = 11; Label = "P*PEA*"; Flag = $false },
    @{ Row = 12; Label = "P*OIL*"; Flag = $false },
    @{ Row = 13; Label = "P*DIS*"; Flag = $false },
    @{ Row = 14; Label = "P*HFO*"; Flag = $false },
    @{ Row = 15; Label = "P*HYD*"; Flag = $true  },
    @{ Row = 16; Label = "P*WIN*"; Flag = $true  },
    @{ Row = 17; Label = "P*SOL*"; Flag = $true  },
    @{ Row = 18; Label = "P*BIO*"; Flag = $true  },
    @{ Row = 19; Label = "P*GEO*"; Flag = $true  },
    @{ Row = 20; Label = "P*OCE*"; Flag = $true  }
)

foreach ($fr in $fuelRows) {
    $r = $fr.Row

    $ws.Range("F$r").Value = $fr.Label

    if ($fr.Flag) {
        $ws.Range("A$r").Value = 1
    } else {
        $ws.Range("A$r").ClearContents()
    }

    $ws.Range("J$r").Value = "ELCC,ELCD"
    $ws.Range("K$r").Value = "UC_FLO"
    $ws.Range("L$r").Value = "O"
    $ws.Range("M$r").Value = "UP"

    $ws.Range("N$r`:T$r").Formula = '=IF($A' + $r + '=1,C$2-1,C$2)'
}

# Row 9 keeps its trailing UC_RES metadata (unchanged by the edit, but make
# sure it is still present).
$ws.Range("C9").Value = "UC_RES"
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = "RES Penetration"
$ws.Range("W9").Value = 15

# ---------------------------------------------------------------
# 4. Apply the new font style (11pt, SimSun, family 3) to the whole
#    label column range, matching the workbook's added cellXfs entry.
# ---------------------------------------------------------------
$ws.Range("F9:F20").Font.Name = "宋体"
$ws.Range("F9:F20").Font.Size = 11
$ws.Range("F9:F20").Font.Family = 3

# ---------------------------------------------------------------
# 5. Selection state (closest reproducible approximation of the
#    recorded activeCell/sqref).
# ---------------------------------------------------------------
$ws.Range("F14:T14").Select()
